$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "51.156.98"
$ws.Range("E2").Value = "  +0.26%  "

$ws.Range("D3").Value = "2.961.28"
$ws.Range("E3").Value = "  +0.94%  "

$ws.Range("E4").Value = "  +0.12%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "379.79"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.59%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "102.37"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +2.03%  "

$ws.Range("E7").Value = "  +2.29%  "

$ws.Range("E8").Value = "  +0.05%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.588"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.53%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "36.48"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.23%  "

$ws.Range("E11").Value = "  -1.07%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0859"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.75%  "

$ws.Range("D13").Value = "3.426.46"
$ws.Range("E13").Value = "  +1.13%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.77"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +4.27%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "18.27"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.78%  "

$ws.Range("D16").Value = "2.969.77"
$ws.Range("E16").Value = "  +1.25%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "11.17"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.47%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.991"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +2.41%  "

$ws.Range("D19").Value = "51.190.17"
$ws.Range("E19").Value = "  +0.38%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "3.16"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.69%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.52"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.94%  "

$ws.Range("E22").Value = "  +0.85%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "70.08"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +2.43%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "266.54"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.97%  "

$ws.Range("E25").Value = "  +2.91%  "

$ws.Range("B26").Value = "Filecoin"
$ws.Range("C26").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.81"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -3.28%  "

$ws.Range("B27").Value = "RenderToken"
$ws.Range("C27").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.44"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +2.48%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.999"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.15%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "25.93"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.66%  "

$ws.Range("E30").Value = "  +1.25%  "

$ws.Range("E31").Value = "  -1.08%  "

$ws.Range("E32").Value = "  +4.49%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "34.55"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +5.13%  "

$ws.Range("E34").Value = "  +1.29%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.02"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.03%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0436"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.43%  "

$ws.Range("E37").Value = "  +0.11%  "

$ws.Range("E38").Value = "  +3.61%  "

$ws.Range("E39").Value = "  +1.55%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "16.60"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +2.38%  "

$ws.Range("E41").Value = "  +3.53%  "

$ws.Range("E42").Value = "  +3.34%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "124.62"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +4.17%  "

$ws.Range("E44").Value = "  +9.84%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "21.43"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +2.08%  "

$ws.Range("B46").Value = "ApeXProtocol"
$ws.Range("C46").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.38"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +3.59%  "

$ws.Range("B47").Value = "WEMIXToken"
$ws.Range("C47").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.02"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.03%  "

$ws.Range("E48").Value = "  -1.51%  "

$ws.Range("D49").Value = "2.034.46"
$ws.Range("E49").Value = "  +2.88%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0327"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.98%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.514"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +12.01%  "

